$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Primordial Hydra', ['{X}{G}{G}', 'Creature — Hydra', 'Primordial Hydra enters the battlefield with X +1/+1 counters on it.', 'At the beginning of your upkeep, double the number of +1/+1 counters on Primordial Hydra.', 'Primordial Hydra has trample as long as it has ten or more +1/+1 counters on it.', '0/0'])"
$ws.Range("A3").Value = "('Serra Avatar', ['{4}{W}{W}{W}', 'Creature — Avatar', 'Serra Avatar’s power and toughness are each equal to your life total.', 'When Serra Avatar is put into a graveyard from anywhere, shuffle it into its owner’s library.', '*/*'])"
$ws.Range("A4").Value = "('Vampire Nocturnus', ['{1}{B}{B}{B}', 'Creature — Vampire', 'Play with the top card of your library revealed.', 'As long as the top card of your library is black, Vampire Nocturnus and other Vampire creatures you control get +2/+1 and have flying.', '3/3'])"

$ws.Range("A5:A20").EntireRow.Delete()
